$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Étude marché + concurrence (Samsung, LG, Bosch)"
$ws.Range("E2").Value = "Product Owner / Business"

# Row 3
$ws.Range("B3").Value = "Définition fonctionnelle du Smart Fridge (détection, seuils, Drive)"
$ws.Range("C3").Value = 1

# Row 4
$ws.Range("A4").Value = "Prototype / Pilote"
$ws.Range("B4").Value = "Prototype IA embarquée (caméra + Jetson/RPi dans un frigo modifié)"
$ws.Range("E4").Value = "IA Engineer / IoT Engineer"

# Row 5
$ws.Range("A5").Value = "Prototype / Pilote"
$ws.Range("B5").Value = "Développement app mobile + backend (inventaire, panier Drive)"
$ws.Range("E5").Value = "Mobile Dev / Backend Dev"

# Row 6
$ws.Range("A6").Value = "Prototype / Pilote"
$ws.Range("B6").Value = "Intégration API Drive (ajout auto au panier Leclerc)"
$ws.Range("E6").Value = "Backend Dev / Business Retail"

# Row 7
$ws.Range("B7").Value = "Intégration hardware propre dans frigo de série"

# Row 8
$ws.Range("B8").Value = "Validation RGPD, conformité CE, sécurité des données"

# Row 9
$ws.Range("B9").Value = "Communication interne équipe (planning, priorités, responsabilités)"
$ws.Range("E9").Value = "Product Owner"

# Row 10
$ws.Range("B10").Value = "Communication Board / Investisseur (demande de budget, vision marché)"
$ws.Range("E10").Value = "Product Owner / Business"

# Row 11
$ws.Range("A11").Value = "Soutenance finale"
$ws.Range("B11").Value = "Prépa soutenance + ajustement livrables (Gantt, budget, risques)"
$ws.Range("E11").Value = "Toute l'équipe"

$wb.Save()
